$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "taxi game"

$ws.Range("A9").Value = "passive income ideas"
$ws.Range("B9").Value = "passive.income.nadi.myfirstdrawermenuproject"
$ws.Rows.Item(9).RowHeight = 24

$ws.Range("A9:B9").Copy()
$ws.Range("A10:B10").PasteSpecial(-4122)
$ws.Range("A10").Value = "passive income apps"
$ws.Range("B10").Value = "passive.income.nadi.myfirstdrawermenuproject"
$ws.Rows.Item(10).RowHeight = 24

$ws.Rows.Item(4).Select()
